# The commit inserts one new price-record row (row 46) into the daily
# "Ají" (chili pepper) price log. All existing rows from 46 downward
# shift down by one to make room; the new row is then populated with
# the latest observation (2022-11-30, 40 units, $15.000 min/max/avg,
# $/caja 15 kilos, Región del Maule, $1.000/kg, 15 kg).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 46, pushing the old row 46 (and all
# rows below it) down to row 47 onward. This also naturally creates the
# new row 112 from the old row 111's content, matching the diff's
# dimension change from A1:R111 to A1:R112.
$ws.Rows("46:46").Insert()

# Populate the newly inserted row 46 with the new observation.
$ws.Range("A46").Value = 7
$ws.Range("B46").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C46").Value = "Ñuble"
$ws.Range("D46").Value = "2022-11-30"
$ws.Range("E46").Value = 16
$ws.Range("F46").Value = 100112021
$ws.Range("G46").Value = "Ají"
$ws.Range("H46").Value = "Americana (o)"
$ws.Range("I46").Value = "Primera"
$ws.Range("J46").Value = 40
$ws.Range("K46").Value = 15000
$ws.Range("L46").Value = 15000
$ws.Range("M46").Value = 15000
$ws.Range("N46").Value = "`$/caja 15 kilos"
$ws.Range("O46").Value = "Región del Maule"
$ws.Range("P46").Value = 1000
$ws.Range("Q46").Value = 15
$ws.Range("R46").Value = "Hortaliza"
